# New weekly price observation for "Terminal La Palmera de La Serena - Albahaca".
# The new record belongs chronologically right after the current row 20
# (2021-11-12), so it is inserted as row 21 and every subsequent row shifts
# down by one (old row 62 becomes the new row 63).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 21..62 down to 22..63, leaving a blank row 21 to fill in.
$ws.Rows.Item(21).Insert()

$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44526
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 100112052
$ws.Range("G21").Value = "Albahaca"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 800
$ws.Range("K21").Value = 3000
$ws.Range("L21").Value = 4000
$ws.Range("M21").Value = 3500
$ws.Range("N21").Value = "`$/paquete"
$ws.Range("O21").Value = "Región de Arica y Parinacota"
$ws.Range("P21").Value = 3500
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
